# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers, styled like the rest of the header ---
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the formatting (border/bold/centered) from an existing header cell
# onto the three new header cells, without touching their values.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-35): every player gets the team's season record ---
$ws.Range("AC2:AC35").Value = 104
$ws.Range("AD2:AD35").Value = 58
$ws.Range("AE2:AE35").Value = 0
